# Apply updated crypto price/volume data per diff (Mon Mar  4 15:57:00 UTC 2024 GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.409.83'
$ws.Range('E2').Value = '  +6.73%  '
$ws.Range('D3').Value = '3.560.74'
$ws.Range('E3').Value = '  +3.91%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '419.47'
$ws.Range('E5').Value = '  +1.34%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '132.07'
$ws.Range('E6').Value = '  +2.11%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.662'
$ws.Range('E7').Value = '  +6.34%  '
$ws.Range('D8').Value = '3.552.42'
$ws.Range('E8').Value = '  +3.86%  '
$ws.Range('E9').Value = '  +0.04%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.787'
$ws.Range('E10').Value = '  +8.74%  '
$ws.Range('E11').Value = '  +20.71%  '
$ws.Range('E12').Value = '  +33.56%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '43.52'
$ws.Range('E13').Value = '  +2.09%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.16'
$ws.Range('E14').Value = '  +9.25%  '
$ws.Range('D15').Value = '4.126.14'
$ws.Range('E15').Value = '  +3.92%  '
$ws.Range('E16').Value = '  +0.14%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '20.51'
$ws.Range('E17').Value = '  -0.15%  '
$ws.Range('D18').Value = '3.540.36'
$ws.Range('E18').Value = '  +3.30%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.12'
$ws.Range('E19').Value = '  +4.76%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.73'
$ws.Range('E20').Value = '  -0.25%  '
$ws.Range('D21').Value = '66.291.28'
$ws.Range('E21').Value = '  +6.45%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '450.37'
$ws.Range('E22').Value = '  -4.45%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '90.33'
$ws.Range('E23').Value = '  -0.72%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.24'
$ws.Range('E24').Value = '  -1.50%  '
$ws.Range('E25').Value = '  -2.19%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.38'
$ws.Range('E26').Value = '  +2.70%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.08'
$ws.Range('E27').Value = '  -4.01%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '34.38'
$ws.Range('E28').Value = '  +3.58%  '
$ws.Range('E29').Value = '  +0.97%  '
$ws.Range('E30').Value = '  +6.41%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '12.52'
$ws.Range('E31').Value = '  +4.65%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.119'
$ws.Range('E32').Value = '  +6.70%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.32'
$ws.Range('E33').Value = '  -4.27%  '
$ws.Range('E34').Value = '  -2.42%  '
$ws.Range('E35').Value = '  -0.03%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '39.09'
$ws.Range('E36').Value = '  -3.69%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '57.57'
$ws.Range('E37').Value = '  -1.98%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0505'
$ws.Range('E38').Value = '  +3.40%  '
$ws.Range('D39').Value = '0.0₃0766'
$ws.Range('E39').Value = '  +42.40%  '
$ws.Range('E40').Value = '  +11.17%  '
$ws.Range('B41').Value = 'FirstDigitalUSD'
$ws.Range('C41').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.999'
$ws.Range('E41').Value = '  -0.06%  '
$ws.Range('B42').Value = 'Stacks'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.06'
$ws.Range('E42').Value = '  +0.77%  '
$ws.Range('E43').Value = '  +4.07%  '
$ws.Range('B44').Value = 'Monero'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '147.97'
$ws.Range('E44').Value = '  +1.95%  '
$ws.Range('B45').Value = 'NEARProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.45'
$ws.Range('E45').Value = '  +3.46%  '
$ws.Range('E46').Value = '  -1.63%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.311'
$ws.Range('E47').Value = '  -3.83%  '
$ws.Range('E48').Value = '  -2.81%  '
$ws.Range('E49').Value = '  -3.37%  '
$ws.Range('E50').Value = '  +6.58%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '15.70'
$ws.Range('E51').Value = '  -3.92%  '
